$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Resize the table's grid columns (w:tblGrid/w:gridCol) to the new widths.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$newColWidthsPt = @(87.7, 125.2, 92.35, 61.3, 86.5)
for ($i = 1; $i -le $t.Columns.Count; $i++) {
    $col = $t.Columns.Item($i)
    $col.Width = $newColWidthsPt[$i - 1]
}

# ---------------------------------------------------------------------------
# 2) Append the new content block at the end of the document (just before
#    the final section break), mirroring the earlier "AND/OR" block above
#    the table.
# ---------------------------------------------------------------------------

# Start a new paragraph after the last (empty) paragraph in the body.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# --- Paragraph: "AND/OR " -------------------------------------------------
$p1 = $word.ActiveDocument.Content
$p1.Collapse(0)
$p1.InsertAfter("AND/OR ")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# --- Empty paragraph --------------------------------------------------------
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# --- Paragraph: "Vanligvis skal ... «3.1.13.txt»." --------------------------
$r1 = $word.ActiveDocument.Content
$r1.Collapse(0)
$r1.InsertAfter("Vanligvis skal alle dokumenter være ferdigstilt, men her finner vi ")

$r2 = $word.ActiveDocument.Content
$r2.Collapse(0)
$r2.InsertAfter("ANTALL")
$r2.Font.Bold = 1
$r2.Font.Underline = 1

$r3 = $word.ActiveDocument.Content
$r3.Collapse(0)
$r3.InsertAfter(" dokumenter som er angitt ")

$r4 = $word.ActiveDocument.Content
$r4.Collapse(0)
$r4.InsertAfter([char]0x00AB)
$r4.Font.Underline = 1

$r5 = $word.ActiveDocument.Content
$r5.Collapse(0)
$r5.InsertAfter("ANNET ENN FERDIGSTILT")
$r5.Font.Bold = 1
$r5.Font.Underline = 1

$r6 = $word.ActiveDocument.Content
$r6.Collapse(0)
$r6.InsertAfter([char]0x00BB)
$r6.Font.Underline = 1

$r7 = $word.ActiveDocument.Content
$r7.Collapse(0)
$r7.InsertAfter(" i journalposter og mapper som ikke utgår")

$r8 = $word.ActiveDocument.Content
$r8.Collapse(0)
$r8.InsertAfter(". Oversikt over disse finnes i vedlegget «3.1.13.txt».")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# --- Empty paragraph --------------------------------------------------------
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# --- Paragraph: "Konsekvensvurdering: Disse filene ..." ---------------------
$k1 = $word.ActiveDocument.Content
$k1.Collapse(0)
$k1.InsertAfter("Konsekvensvurdering: ")
$k1.Font.Bold = 1

$k2 = $word.ActiveDocument.Content
$k2.Collapse(0)
$k2.InsertAfter("Disse filene kan ikke anerkjennes som fullstendige filer.")

Write-Host "Edit complete"
